$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" column (C) for all data rows from 46079 to 46081
$ws.Range("C2").Value = 46081
$ws.Range("C3").Value = 46081
$ws.Range("C4").Value = 46081
$ws.Range("C5").Value = 46081
$ws.Range("C6").Value = 46081
$ws.Range("C7").Value = 46081
$ws.Range("C8").Value = 46081

# Swap the data between row 6 and row 7 (columns A, B, G)
$a6 = $ws.Range("A6").Value()
$b6 = $ws.Range("B6").Value()
$g6 = $ws.Range("G6").Value()

$a7 = $ws.Range("A7").Value()
$b7 = $ws.Range("B7").Value()
$g7 = $ws.Range("G7").Value()

$ws.Range("A6").Value = $a7
$ws.Range("B6").Value = $b7
$ws.Range("G6").Value = $g7

$ws.Range("A7").Value = $a6
$ws.Range("B7").Value = $b6
$ws.Range("G7").Value = $g6
